# Updates the "cryptos" worksheet with refreshed price/volume-change data
# (commit: "Updated cryptos list ... with GitHub Actions").
#
# For most rows only the Price (D) and/or Volume(1h) (E) columns change.
# Rows 16/17, 44/45 and 47/48 also swap rank order, so the Coin (B) and
# Link (C) columns are rewritten too so the two coins trade places.
#
# A couple of Price values (rows 29 and 31) are exact multiples of 0.10,
# which Excel would otherwise normalise to "9.7"/"2.8" when it auto-detects
# them as numbers. Those two cells are briefly switched to a text number
# format before the assignment (and switched back to the default "Normal"
# style afterwards) so the literal trailing-zero text such as "9.70" is
# preserved, matching how every other cell in this sheet is already stored
# as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.890.43'
$ws.Range("E2").Value = '  +4.38%  '

$ws.Range("D3").Value = '2.260.12'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").Value = '301.58'
$ws.Range("E5").Value = '  +0.40%  '

$ws.Range("D6").Value = '100.22'
$ws.Range("E6").Value = '  +6.48%  '

$ws.Range("D7").Value = '0.559'
$ws.Range("E7").Value = '  -1.20%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = '0.505'
$ws.Range("E9").Value = '  -0.49%  '

$ws.Range("E10").Value = '  +3.82%  '

$ws.Range("D11").Value = '0.0773'
$ws.Range("E11").Value = '  -2.16%  '

$ws.Range("D12").Value = '7.08'
$ws.Range("E12").Value = '  -1.58%  '

$ws.Range("E13").Value = '  -1.23%  '

$ws.Range("D14").Value = '2.603.33'
$ws.Range("E14").Value = '  -0.35%  '

$ws.Range("D15").Value = '2.255.36'
$ws.Range("E15").Value = '  -0.86%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '13.54'
$ws.Range("E16").Value = '  -0.29%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '46.879.08'
$ws.Range("E17").Value = '  +4.56%  '

$ws.Range("D18").Value = '0.789'
$ws.Range("E18").Value = '  -1.42%  '

$ws.Range("D19").Value = '12.67'
$ws.Range("E19").Value = '  -5.19%  '

$ws.Range("D20").Value = '0.0₃0929'
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("D21").Value = '5.79'
$ws.Range("E21").Value = '  -3.85%  '

$ws.Range("D22").Value = '65.19'
$ws.Range("E22").Value = '  -0.38%  '

$ws.Range("D23").Value = '247.61'
$ws.Range("E23").Value = '  +3.57%  '

$ws.Range("E24").Value = '  -1.99%  '

$ws.Range("E25").Value = '  +0.26%  '

$ws.Range("D26").Value = '1.86'
$ws.Range("E26").Value = '  -2.22%  '

$ws.Range("D27").Value = '42.08'
$ws.Range("E27").Value = '  +2.04%  '

$ws.Range("E28").Value = '  +0.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.70'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.55%  '

$ws.Range("D30").Value = '20.02'
$ws.Range("E30").Value = '  +2.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +9.77%  '

$ws.Range("D32").Value = '145.76'
$ws.Range("E32").Value = '  -4.45%  '

$ws.Range("D33").Value = '5.35'
$ws.Range("E33").Value = '  -2.79%  '

$ws.Range("D34").Value = '3.24'
$ws.Range("E34").Value = '  +11.62%  '

$ws.Range("D35").Value = '0.0764'
$ws.Range("E35").Value = '  -2.89%  '

$ws.Range("E36").Value = '  +11.69%  '

$ws.Range("E37").Value = '  -1.88%  '

$ws.Range("D38").Value = '16.09'
$ws.Range("E38").Value = '  +18.51%  '

$ws.Range("D39").Value = '1.68'
$ws.Range("E39").Value = '  -4.36%  '

$ws.Range("D40").Value = '3.86'
$ws.Range("E40").Value = '  -0.34%  '

$ws.Range("D41").Value = '0.0296'
$ws.Range("E41").Value = '  -3.42%  '

$ws.Range("D42").Value = '3.12'
$ws.Range("E42").Value = '  -2.33%  '

$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.13%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.92'
$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("B45").Value = 'BitcoinSV'
$ws.Range("C45").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D45").Value = '91.44'
$ws.Range("E45").Value = '  +19.55%  '

$ws.Range("D46").Value = '1.768.23'
$ws.Range("E46").Value = '  -0.37%  '

$ws.Range("B47").Value = 'ordi'
$ws.Range("C47").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D47").Value = '71.14'
$ws.Range("E47").Value = '  +2.56%  '

$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.184'
$ws.Range("E48").Value = '  -4.04%  '

$ws.Range("D49").Value = '4.81'
$ws.Range("E49").Value = '  +2.73%  '

$ws.Range("D50").Value = '7.82'
$ws.Range("E50").Value = '  -0.67%  '

$ws.Range("D51").Value = '93.42'
$ws.Range("E51").Value = '  -2.03%  '

